$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Move the existing totals row (row 37: "Yht" / SUM / empty) down to row 38,
# preserving its style, then rebuild its formula on the new row.
# ---------------------------------------------------------------------
$ws.Range("B37:D37").Copy()
$ws.Range("B38:D38").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(38, 2).Value = "Yht"
$ws.Cells.Item(38, 3).Formula = "=SUM(C6:C37)"
$ws.Cells.Item(38, 4).Value = ""
$ws.Rows.Item(38).RowHeight = 18.75

# ---------------------------------------------------------------------
# Turn (old) row 37 into a new time-entry row, copying the style of the
# row above it (row 36) and filling in the new entry's data.
# ---------------------------------------------------------------------
$ws.Range("B36:D36").Copy()
$ws.Range("B37:D37").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(37, 2).Value = 45360
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(37, 4).Value = "Tein viimeistelyä projektiin sen tiedostoihin ja rakenteeseen ja palautin projektin"
$ws.Rows.Item(37).RowHeight = 37.5

# ---------------------------------------------------------------------
# Update the sheet view / selection to match.
# ---------------------------------------------------------------------
$ws.Range("B41").Select()
